# Add two new worksheets "rho_min" and "rho_max" after "rho_mat",
# each containing the min/max purity coherency-matrix cases computed
# by computeRho.m, with the same layout/header style as "rho_mat".

$wb = $excel.ActiveWorkbook

# ---- small full-precision touch-ups to the existing rho_mat values ----
$rhoMat = $wb.Worksheets.Item("rho_mat")
$rhoMat.Range("B4").Value  = 0.197125134359208
$rhoMat.Range("C4").Value  = 0.80287486564079191
$rhoMat.Range("D4").Value  = -0.26262667603843348
$rhoMat.Range("F4").Value  = 0.99999999999999978
$rhoMat.Range("B5").Value  = 0.36971967045748444
$rhoMat.Range("D5").Value  = -0.23567565231848192
$rhoMat.Range("D6").Value  = -0.11099791392217158
$rhoMat.Range("B12").Value = 0.065348071014497328
$rhoMat.Range("B13").Value = 0.21487636310143046
$rhoMat.Range("E13").Value = 0.3001146734777696
$rhoMat.Range("D15").Value = -0.10972222409935729

# ---- header labels/style, reused verbatim for both new sheets ----
$headers = New-Object 'object[,]' 1,6
$headers[0,0] = "theta"
$headers[0,1] = "Jxx"
$headers[0,2] = "Jyy"
$headers[0,3] = "beta"
$headers[0,4] = "gamma"
$headers[0,5] = "trace_sq"

# ---- rho_min : minimum-purity case ----
$data4 = New-Object 'object[,]' 19,6
$data4[0,0] = 0
$data4[0,1] = 0.10175205709792243
$data4[0,2] = 0.89824794290207766
$data4[0,3] = -0.0065757803645066174
$data4[0,4] = 0.027449599740522294
$data4[0,5] = 1
$data4[1,0] = 10
$data4[1,1] = 0.13798760510086441
$data4[1,2] = 0.86201239489913561
$data4[1,3] = -0.13713989959826944
$data4[1,4] = -0.11151274664303727
$data4[1,5] = 1
$data4[2,0] = 20
$data4[2,1] = 0.25544078274679349
$data4[2,2] = 0.74455921725320662
$data4[2,3] = -0.21206043025684346
$data4[2,4] = -0.24128595813127063
$data4[2,5] = 1
$data4[3,0] = 30
$data4[3,1] = 0.39517371228991705
$data4[3,2] = 0.60482628771008307
$data4[3,3] = -0.1896295766440827
$data4[3,4] = -0.33897759104451547
$data4[3,5] = 1
$data4[4,0] = 40
$data4[4,1] = 0.48739018448041027
$data4[4,2] = 0.51260981551958962
$data4[4,3] = -0.089338879404745442
$data4[4,4] = -0.39219057139218494
$data4[4,5] = 0.99999999999999978
$data4[5,0] = 50
$data4[5,1] = 0.49257577532103869
$data4[5,2] = 0.50742422467896131
$data4[5,3] = 0.044146846062254011
$data4[5,4] = -0.39826593212311523
$data4[5,5] = 1
$data4[6,0] = 60
$data4[6,1] = 0.40375077069942877
$data4[6,2] = 0.59624922930057123
$data4[6,3] = 0.14762864880643167
$data4[6,4] = -0.35454124319654845
$data4[6,5] = 1
$data4[7,0] = 70
$data4[7,1] = 0.26878299797580446
$data4[7,2] = 0.73121700202419548
$data4[7,3] = 0.17110979333065468
$data4[7,4] = -0.27151877030530536
$data4[7,5] = 1
$data4[8,0] = 80
$data4[8,1] = 0.15126976246911178
$data4[8,2] = 0.84873023753088828
$data4[8,3] = 0.10471132793586194
$data4[8,4] = -0.15914276712719105
$data4[8,5] = 1
$data4[9,0] = 90
$data4[9,1] = 0.10195184483067134
$data4[9,2] = 0.89804815516932868
$data4[9,3] = -0.027298730454168488
$data4[9,4] = -0.021722089799595144
$data4[9,5] = 1
$data4[10,0] = 100
$data4[10,1] = 0.14869508965833422
$data4[10,2] = 0.85130491034166578
$data4[10,3] = -0.15985552296083849
$data4[10,4] = 0.11977375466682297
$data4[10,5] = 1
$data4[11,0] = 110
$data4[11,1] = 0.26781780891151341
$data4[11,2] = 0.73218219108848659
$data4[11,3] = -0.22835238175094452
$data4[11,4] = 0.24438970835190182
$data4[11,5] = 1
$data4[12,0] = 120
$data4[12,1] = 0.40999229383721336
$data4[12,2] = 0.59000770616278664
$data4[12,3] = -0.20149178189256522
$data4[12,4] = 0.3414501453628499
$data4[12,5] = 1
$data4[13,0] = 130
$data4[13,1] = 0.5040503136530351
$data4[13,2] = 0.4959496863469649
$data4[13,3] = -0.088832196454909459
$data4[13,4] = 0.39491704728323862
$data4[13,5] = 1
$data4[14,0] = 140
$data4[14,1] = 0.50410689069681258
$data4[14,2] = 0.49589310930318742
$data4[14,3] = 0.046448830827355569
$data4[14,4] = 0.39768818558885366
$data4[14,5] = 1
$data4[15,0] = 150
$data4[15,1] = 0.41452249512373757
$data4[15,2] = 0.58547750487626238
$data4[15,3] = 0.15201262640044458
$data4[15,4] = 0.35303866972439696
$data4[15,5] = 1
$data4[16,0] = 160
$data4[16,1] = 0.27829056316663581
$data4[16,2] = 0.72170943683336419
$data4[16,3] = 0.17859272382538791
$data4[16,4] = 0.26977396936684583
$data4[16,5] = 1
$data4[17,0] = 170
$data4[17,1] = 0.15241825318350369
$data4[17,2] = 0.84758174681649634
$data4[17,3] = 0.11064804073844839
$data4[17,4] = 0.15163757188995072
$data4[17,5] = 1
$data4[18,0] = 180
$data4[18,1] = 0.10160249830240942
$data4[18,2] = 0.89839750169759058
$data4[18,3] = -0.012872261587301235
$data4[18,4] = 0.021251777309605303
$data4[18,5] = 1

$rhoMin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rhoMat)
$rhoMin.Name = "rho_min"
$hdrRange = $rhoMin.Range("A1:F1")
$hdrRange.NumberFormat = "@"
$hdrRange.Value = $headers
$rhoMin.Range("A2:F20").Value = $data4

# ---- rho_max : maximum-purity case ----
$data5 = New-Object 'object[,]' 19,6
$data5[0,0] = 0
$data5[0,1] = 0.0012511435290904109
$data5[0,2] = 0.99874885647090961
$data5[0,3] = -0.0082352288207749764
$data5[0,4] = 0.03437671612665661
$data5[0,5] = 1
$data5[1,0] = 10
$data5[1,1] = 0.050695970049746394
$data5[1,2] = 0.94930402995025354
$data5[1,3] = -0.17020828685615477
$data5[1,4] = -0.13840168779717638
$data5[1,5] = 0.99999999999999978
$data5[2,0] = 20
$data5[2,1] = 0.19712513464556097
$data5[2,2] = 0.80287486535443897
$data5[2,3] = -0.26262667579013027
$data5[2,4] = -0.29882109086594749
$data5[2,5] = 1
$data5[3,0] = 30
$data5[3,1] = 0.36971967056557414
$data5[3,2] = 0.6302803294344258
$data5[3,3] = -0.23567565212295022
$data5[3,4] = -0.42128852596885136
$data5[3,5] = 1
$data5[4,0] = 40
$data5[4,1] = 0.4843331008158922
$data5[4,2] = 0.51566689918410769
$data5[4,3] = -0.11099791386169733
$data5[4,4] = -0.48727200912761037
$data5[4,5] = 0.99999999999999978
$data5[5,0] = 50
$data5[5,1] = 0.49073764230674005
$data5[5,2] = 0.50926235769326
$data5[5,3] = 0.055076980686835829
$data5[5,4] = -0.4968709434154655
$data5[5,5] = 1
$data5[6,0] = 60
$data5[6,1] = 0.37845056996600768
$data5[6,2] = 0.62154943003399232
$data5[6,3] = 0.1864345122502056
$data5[6,4] = -0.44773642705758043
$data5[6,5] = 1
$data5[7,0] = 70
$data5[7,1] = 0.20772955137034738
$data5[7,2] = 0.79227044862965257
$data5[7,3] = 0.21629177622692419
$data5[7,4] = -0.34321400292266979
$data5[7,5] = 1
$data5[8,0] = 80
$data5[8,1] = 0.061203141469662152
$data5[8,2] = 0.9387968585303379
$data5[8,3] = 0.1317551413841094
$data5[8,4] = -0.20024459813884465
$data5[8,5] = 1
$data5[9,0] = 90
$data5[9,1] = 0.001909374715063373
$data5[9,2] = 0.99809062528493664
$data5[9,3] = -0.034159790831381931
$data5[9,4] = -0.027181558688981422
$data5[9,5] = 1
$data5[10,0] = 100
$data5[10,1] = 0.065348071297809132
$data5[10,2] = 0.93465192870219094
$data5[10,3] = -0.19778121319468822
$data5[10,4] = 0.1481900535441027
$data5[10,5] = 1
$data5[11,0] = 110
$data5[11,1] = 0.21487636340786598
$data5[11,2] = 0.78512363659213402
$data5[11,3] = -0.28042056629783085
$data5[11,4] = 0.30011467315522566
$data5[11,5] = 1
$data5[12,0] = 120
$data5[12,1] = 0.38930491379431426
$data5[12,2] = 0.6106950862056858
$data5[12,3] = -0.247802672873318
$data5[12,4] = 0.41992908037813764
$data5[12,5] = 1
$data5[13,0] = 130
$data5[13,1] = 0.50500279673189175
$data5[13,2] = 0.49499720326810831
$data5[13,3] = -0.10972222405007093
$data5[13,4] = 0.48778684387477178
$data5[13,5] = 1
$data5[14,0] = 140
$data5[14,1] = 0.50512832344926251
$data5[14,2] = 0.49487167655073755
$data5[14,3] = 0.058001209651775895
$data5[14,4] = 0.49659798572989566
$data5[14,5] = 1
$data5[15,0] = 150
$data5[15,1] = 0.3914611559334199
$data5[15,2] = 0.60853884406657999
$data5[15,3] = 0.19302475870012303
$data5[15,4] = 0.44828647230822971
$data5[15,5] = 0.99999999999999978
$data5[16,0] = 160
$data5[16,1] = 0.21735853277377476
$data5[16,2] = 0.78264146722622518
$data5[16,3] = 0.22767506074121016
$data5[16,4] = 0.34391549412755357
$data5[16,5] = 1
$data5[17,0] = 170
$data5[17,1] = 0.060058368611865906
$data5[17,2] = 0.93994163138813402
$data5[17,3] = 0.14004958545211904
$data5[17,4] = 0.19193090939904967
$data5[17,5] = 0.99999999999999978
$data5[18,0] = 180
$data5[18,1] = 0.00096953174798121867
$data5[18,2] = 0.9990304682520188
$data5[18,3] = -0.016123722412921596
$data5[18,4] = 0.026619856642700692
$data5[18,5] = 1

$rhoMax = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rhoMin)
$rhoMax.Name = "rho_max"
$hdrRange2 = $rhoMax.Range("A1:F1")
$hdrRange2.NumberFormat = "@"
$hdrRange2.Value = $headers
$rhoMax.Range("A2:F20").Value = $data5

Write-Host "Added rho_min and rho_max sheets"
